$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Inside the empty function, add this:" paragraph: numId 5 -> 4
# ---------------------------------------------------------------------------
$pNum = $d.Paragraphs.Item(19)
if ($pNum.Range.Text.Trim() -ne "Inside the empty function, add this:") {
    throw "Paragraph 19 text mismatch: $($pNum.Range.Text)"
}
$xml = $pNum.Range.WordOpenXML
$startTag = $xml.IndexOf("<w:p ")
if ($startTag -lt 0) { $startTag = $xml.IndexOf("<w:p>") }
$endTag = $xml.IndexOf("</w:p>") + 6
$frag = $xml.Substring($startTag, $endTag - $startTag)
$frag2 = $frag.Replace('<w:numId w:val="5"/>', '<w:numId w:val="4"/>')
if ($frag2 -eq $frag) { throw "numId 5 not found in paragraph 19 fragment" }
$pNum.Range.InsertXML($frag2)

# ---------------------------------------------------------------------------
# 2) "  bodytext: ..." Heading3 paragraph: add <w:ind w:left="810" w:firstLine="0"/>
# ---------------------------------------------------------------------------
$pInd = $d.Paragraphs.Item(21)
if (-not $pInd.Range.Text.Contains("bodytext")) {
    throw "Paragraph 21 text mismatch: $($pInd.Range.Text)"
}
$pInd.Range.ParagraphFormat.LeftIndent = 40.5
$pInd.Range.ParagraphFormat.FirstLineIndent = 0

# ---------------------------------------------------------------------------
# 3) "Use req.body.[element's name] ..." -> "Use req.body.[[element's name]] ..."
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute("req.body.[element's name]", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "req.body.[[element's name]]", 2)
if (-not $found) { throw "Find/replace for req.body text failed" }

# ---------------------------------------------------------------------------
# 4) Swap the docPr/cNvPr "name" labels between the two screenshots
#    (image1.png <-> image2.png) while keeping each run's real embedded
#    picture (r:embed) untouched.
# ---------------------------------------------------------------------------
$pImg1 = $d.Paragraphs.Item(26)
$xml1 = $pImg1.Range.WordOpenXML
$s1 = $xml1.IndexOf("<w:p ")
if ($s1 -lt 0) { $s1 = $xml1.IndexOf("<w:p>") }
$e1 = $xml1.IndexOf("</w:p>") + 6
$frag1 = $xml1.Substring($s1, $e1 - $s1)
if (-not $frag1.Contains('name="image1.png"')) { throw "Paragraph 26 does not contain image1.png" }
$frag1b = $frag1.Replace('name="image1.png"', 'name="image2.png"')
$frag1b = $frag1b.Replace('r:embed="rId5"', 'r:embed="rId6"')
$pImg1.Range.InsertXML($frag1b)

$pImg2 = $d.Paragraphs.Item(28)
$xml2 = $pImg2.Range.WordOpenXML
$s2 = $xml2.IndexOf("<w:p ")
if ($s2 -lt 0) { $s2 = $xml2.IndexOf("<w:p>") }
$e2 = $xml2.IndexOf("</w:p>") + 6
$frag2x = $xml2.Substring($s2, $e2 - $s2)
if (-not $frag2x.Contains('name="image2.png"')) { throw "Paragraph 28 does not contain image2.png" }
$frag2b = $frag2x.Replace('name="image2.png"', 'name="image1.png"')
$frag2b = $frag2b.Replace('r:embed="rId5"', 'r:embed="rId7"')
$pImg2.Range.InsertXML($frag2b)

# ---------------------------------------------------------------------------
# 5) Delete the whole "Converting POST data to JSON" section (9 paragraphs),
#    right before the section break.
# ---------------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(29)
if (-not $pStart.Range.Text.Contains("Converting POST data to JSON")) {
    throw "Paragraph 29 text mismatch: $($pStart.Range.Text)"
}
$pEnd = $d.Paragraphs.Item(37)
if (-not $pEnd.Range.Text.Contains("The POST data is now stored as a JSON object")) {
    throw "Paragraph 37 text mismatch: $($pEnd.Range.Text)"
}
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

Write-Host "Done. Paragraph count now:" $d.Paragraphs.Count
